# The "Metadata" sheet's header row used TitleCase/PascalCase column names;
# normalize them all to lowercase (values themselves change, the columns
# they sit in stay the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "sample_id"
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "lifestage"
$ws.Range("D1").Value = "date_collected"
$ws.Range("E1").Value = "experiment"
$ws.Range("F1").Value = "preservation"
$ws.Range("G1").Value = "freezer"
$ws.Range("H1").Value = "transport_method"
$ws.Range("I1").Value = "storage_location"

# Scroll the view one column to the right (column B becomes the leftmost
# visible column) and move the active selection to J1, just past the
# last header cell.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J1").Select()
